$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WMS Location Floor Location")

# Row 147: BUBBLEUP value changes from "3" to "30".
# Force text (not numeric) storage for the digit-only value, then restore
# the default "Normal" cell style so formatting is unaffected.
$ws.Range("B147").NumberFormat = "@"
$ws.Range("B147").Value = "30"
$ws.Range("B147").Style = "Normal"

# New row 148: WD40 / 40 / N
$ws.Range("A148").Value = "WD40"
$ws.Range("B148").NumberFormat = "@"
$ws.Range("B148").Value = "40"
$ws.Range("B148").Style = "Normal"
$ws.Range("C148").Value = "N"

# New row 149: ROCK / PAPER / Y
$ws.Range("A149").Value = "ROCK"
$ws.Range("B149").Value = "PAPER"
$ws.Range("C149").Value = "Y"
